$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows whose B (M2 local currency) and D (M2 USD) values changed ---
# Row -> @(B, D)
$updates = @{
    368 = @(1486571000000, 404734886126.9552)
    372 = @(1496040000000, 407301833626.0056)
    379 = @(1485900000000, 404555503523.0825)
    380 = @(1498300000000, 407931563987.2363)
    381 = @(1517100000000, 413038299494.2861)
    382 = @(1563100000000, 425567829892.662)
    383 = @(1563400000000, 425647189827.3752)
    384 = @(1570300000000, 427524601790.632)
    385 = @(1593400000000, 433813730174.6118)
    386 = @(1567300000000, 426707831870.6346)
    387 = @(1568100000000, 426923893689.6613)
    388 = @(1622200000000, 441656525065.41)
    389 = @(1606000000000, 437240592588.9625)
    390 = @(1627700000000, 443153942700.6336)
    391 = @(1645600000000, 448021257595.8879)
    392 = @(1629400000000, 443614969214.5804)
    393 = @(1678100000000, 456875111455.3654)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
}

# --- Append new rows 400-402 (same pattern/format as existing data rows) ---
$newRows = @(
    @(400, 45078, 1855306000000, 0.2722551575336405, 505116627303.1085),
    @(401, 45108, 1858844000000, 0.2722521926510966, 506074354796.335),
    @(402, 45139, 1860330000000, 0.272253304474483, 506480989913.015)
)

foreach ($r in $newRows) {
    $rowIdx = $r[0]

    # Copy formatting from the row directly above so style (e.g. date-ish
    # numFmt on column A) carries forward to the newly appended row.
    $ws.Cells.Item($rowIdx - 1, 1).Copy()
    $ws.Cells.Item($rowIdx, 1).PasteSpecial(-4122)

    $ws.Cells.Item($rowIdx, 1).Value = $r[1]
    $ws.Cells.Item($rowIdx, 2).Value = $r[2]
    $ws.Cells.Item($rowIdx, 3).Value = $r[3]
    $ws.Cells.Item($rowIdx, 4).Value = $r[4]
}
